# Add the "Alpha Diversity" column (column I) with the computed diversity
# values for each of the 25 survey plots, and leave the selection on J5
# (next empty cell to the right of the new column), matching the author's
# "Analyzed data and made graphs" session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Cells.Item(1, 9).Value = "Alpha Diversity"

# Data rows 2-26 (Site numbers 1-25)
$alphaDiversity = @{
    2  = 5
    3  = 4
    4  = 6
    5  = 13
    6  = 5
    7  = 12
    8  = 11
    9  = 11
    10 = 10
    11 = 9
    12 = 12
    13 = 3
    14 = 11
    15 = 14
    16 = 5
    17 = 5
    18 = 3
    19 = 3
    20 = 5
    21 = 4
    22 = 4
    23 = 5
    24 = 6
    25 = 13
    26 = 9
}

foreach ($row in $alphaDiversity.Keys | Sort-Object) {
    $ws.Cells.Item($row, 9).Value = $alphaDiversity[$row]
}

# Move the active selection to J5, mirroring where the author left off.
$ws.Range("J5").Select()
